# Replace every occurrence of the value "congenital" with "misc_long_term"
# across all worksheets in the workbook (each matching sheet has exactly
# one such cell, either A3 or A4 depending on the sheet's layout).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $found = $ws.Cells.Find("congenital")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            if ($found.Text -eq "congenital") {
                $found.Value = "misc_long_term"
            }
            $found = $ws.Cells.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}
